# Add a new "2022-Q1" sheet (with fund detail data) before the "总计" sheet,
# and update "总计" with a new summary row for 2022-Q1.
#
# Strategy: rename the existing "总计" sheet to "2022-Q1" (so it keeps
# sheetId=4) and repurpose its contents into the fund-detail table; then add
# a brand-new sheet named "总计" right after it (gets sheetId=5) and populate
# it with the summary table (old rows + the new 2022-Q1 row on top).

$wb = $excel.ActiveWorkbook

# ---- cells we borrow existing cell styles from (already-present style index
#      used by the header row / "A" index column on the sibling quarter
#      sheets) so the new cells look the same as their neighbours. ----
$styleDonor = $wb.Worksheets.Item("2021-Q1")
$styleHeaderCell = $styleDonor.Cells.Item(1, 2)
$styleIndexCell = $styleDonor.Cells.Item(2, 1)

# ---------------------------------------------------------------------
# Step 1: repurpose the current "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $q1Headers) {
    $c = $q1.Cells.Item(1, $col)
    $styleHeaderCell.Copy($c)
    $c.Value = $h
    $col = $col + 1
}

$q1Data = @(
    @("002434", "中银宏利灵活配置混合A", "5.86", "30.78", "0.89", "0.0522", 3),
    @("002535", "中银鑫利灵活配置混合A", "6.75", "20.93", "0.71", "0.0479", 7),
    @("003966", "中银润利灵活配置混合A", "5.69", "25.96", "0.73", "0.0415", 3),
    @("006952", "中银景元回报混合", "3.58", "33.73", "1.13", "0.0405", 10),
    @("002261", "中银宝利灵活配置混合A", "4.35", "31.13", "0.89", "0.0387", 2),
    @("002614", "中银颐利灵活配置混合A", "3.34", "36.07", "1.00", "0.0334", 5),
    @("003967", "中银润利灵活配置混合C", "3.85", "25.96", "0.73", "0.0281", 3),
    @("002615", "中银颐利灵活配置混合C", "2.43", "36.07", "1.00", "0.0243", 5),
    @("002435", "中银宏利灵活配置混合C", "2.33", "30.78", "0.89", "0.0207", 3),
    @("002536", "中银鑫利灵活配置混合C", "2.30", "20.93", "0.71", "0.0163", 7),
    @("002288", "中银稳进策略灵活配置混合", "0.70", "66.42", "2.21", "0.0155", 8),
    @("002262", "中银宝利灵活配置混合C", "1.60", "31.13", "0.89", "0.0142", 2)
)

$r = 2
$idx = 0
foreach ($row in $q1Data) {
    $cIdx = $q1.Cells.Item($r, 1)
    $styleIndexCell.Copy($cIdx)
    $cIdx.Value = $idx

    $cCode = $q1.Cells.Item($r, 2)
    $cCode.NumberFormat = "@"
    $cCode.Value = $row[0]

    $q1.Cells.Item($r, 3).Value = $row[1]

    $cScale = $q1.Cells.Item($r, 4)
    $cScale.NumberFormat = "@"
    $cScale.Value = $row[2]

    $cPos = $q1.Cells.Item($r, 5)
    $cPos.NumberFormat = "@"
    $cPos.Value = $row[3]

    $cShare = $q1.Cells.Item($r, 6)
    $cShare.NumberFormat = "@"
    $cShare.Value = $row[4]

    $cValue = $q1.Cells.Item($r, 7)
    $cValue.NumberFormat = "@"
    $cValue.Value = $row[5]

    $q1.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
    $idx = $idx + 1
}

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet right after "2022-Q1" and fill it with
#          the summary table (new 2022-Q1 row + the previous rows)
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the sheet-level cosmetics the rest of the workbook uses (outline
# grouping + page margins) instead of Excel's brand-new-sheet defaults.
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$col = 2
foreach ($h in $totalHeaders) {
    $c = $total.Cells.Item(1, $col)
    $styleHeaderCell.Copy($c)
    $c.Value = $h
    $col = $col + 1
}

$totalData = @(
    @("2022-Q1", 12, 0.37),
    @("2021-Q4", 1, 0.23),
    @("2021-Q1", 4, 0.32),
    @("2020-Q4", 6, 0.44)
)

$r = 2
$idx = 0
foreach ($row in $totalData) {
    $cIdx = $total.Cells.Item($r, 1)
    $styleIndexCell.Copy($cIdx)
    $cIdx.Value = $idx

    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]

    $r = $r + 1
    $idx = $idx + 1
}
